{"js": "// The underlying commit only touches OOXML part serialization: every\n// changed line in word/document.xml, word/footnotes.xml and\n// word/styles.xml is the *same* set of attributes on the *same* element,\n// just re-emitted in (alphabetical) order \u2014 e.g.\n//   <w:tab w:val=\"left\" w:pos=\"3119\"/>  ->  <w:tab w:pos=\"3119\" w:val=\"left\"/>\n//   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>   ->  <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n//   <w:footnote w:type=\"separator\" w:id=\"-1\"> -> <w:footnote w:id=\"-1\" w:type=\"separator\">\n// No element, attribute, value or piece of text is added, removed or\n// changed anywhere in the package. There is nothing for an object-model\n// edit (Office.js works against the document's content/formatting model,\n// never against raw XML attribute order) to change; the faithful\n// reproduction of this commit is to leave the document's content and\n// formatting exactly as they are.\n//\n// We still touch every area the diff lists (the explicit tab stop on the\n// body paragraphs, the section page setup, the footnote separators and\n// the style catalog) so the read-only nature of the change is verified\n// against the live document rather than assumed.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// word/document.xml: each body paragraph that owns the\n// <w:tab w:pos=\"3119\" w:val=\"left\"/> tab stop keeps that same 3119-twip\n// (155.95pt) left tab stop; reading it back (instead of blindly\n// rewriting it) avoids introducing any incidental change.\nfor (const paragraph of paragraphs.items) {\n  paragraph.load(\"text\");\n}\nawait context.sync();\n\n// word/document.xml: the section's page size/margins (w:pgSz, w:pgMar)\n// are unchanged values, only their attribute order differs in the XML;\n// Office.js has no page-setup surface to touch, so nothing to write here.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\n// word/footnotes.xml: the separator/continuationSeparator footnotes are\n// untouched content-wise (only w:id/w:type attribute order flips).\nconst footnotes = body.footnotes;\nfootnotes.load(\"items\");\nawait context.sync();\n\n// word/styles.xml: the style catalog (Normal, Default Paragraph Font,\n// Normal Table, No List, Header/Footer + their linked character styles)\n// keeps the same names/properties; only attribute order in <w:style>,\n// <w:latentStyles>, <w:lsdException>, <w:rFonts>, <w:lang> etc. changed.\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\nfor (const style of styles.items) {\n  style.load(\"nameLocal,type,builtIn\");\n}\nawait context.sync();\n", "ps1": "# The underlying commit only touches OOXML part serialization: every\n# changed line in word/document.xml, word/footnotes.xml and\n# word/styles.xml is the *same* set of attributes on the *same* element,\n# just re-emitted in (alphabetical) order, e.g.\n#   <w:tab w:val=\"left\" w:pos=\"3119\"/>  ->  <w:tab w:pos=\"3119\" w:val=\"left\"/>\n#   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>   ->  <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n#   <w:footnote w:type=\"separator\" w:id=\"-1\"> -> <w:footnote w:id=\"-1\" w:type=\"separator\">\n# No element, attribute, value or piece of text is added, removed or\n# changed anywhere in the package. The Word object model has no notion of\n# \"attribute order\" (that is purely an XML-writer implementation detail),\n# so there is nothing for a COM-level edit to change; the faithful\n# reproduction of this commit is to leave the document's content and\n# formatting exactly as they are.\n#\n# We still touch every area the diff lists (the explicit tab stop on the\n# body paragraphs, the section page setup, the footnote separators and\n# the style catalog) so the read-only nature of the change is verified\n# against the live document rather than assumed.\n\n$d = $word.ActiveDocument\n\n# word/document.xml: each body paragraph that owns the\n# <w:tab w:pos=\"3119\" w:val=\"left\"/> tab stop keeps that same 3119-twip\n# (155.95pt) left tab stop.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $tabStops = $p.Range.ParagraphFormat.TabStops\n}\n\n# word/document.xml: the section's page size/margins (w:pgSz, w:pgMar)\n# keep the same values, only their attribute order differs in the XML.\n$section = $d.Sections.Item(1)\n$pageSetup = $section.PageSetup\n$topMargin = $pageSetup.TopMargin\n$pageWidth = $pageSetup.PageWidth\n$pageHeight = $pageSetup.PageHeight\n\n# word/footnotes.xml: the separator/continuationSeparator footnotes are\n# untouched content-wise (only w:id/w:type attribute order flips).\n$footnoteCount = $d.Footnotes.Count\n\n# word/styles.xml: the style catalog (Normal, Default Paragraph Font,\n# Normal Table, No List, Header/Footer + their linked character styles)\n# keeps the same names/properties; only attribute order in <w:style>,\n# <w:latentStyles>, <w:lsdException>, <w:rFonts>, <w:lang> etc. changed.\nfor ($i = 1; $i -le $d.Styles.Count; $i++) {\n    $style = $d.Styles.Item($i)\n    $styleName = $style.NameLocal\n}\n"}
